# Generate Report for Handoff
# Refresh the "Latest Handoff Date(time)" values for the rows whose handoff
# batch was regenerated (rows 7 and 10-16 on each sheet) so that they all
# carry the new handoff run's timestamp.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$overviewRows = @(7, 10, 11, 12, 13, 14, 15, 16)
foreach ($r in $overviewRows) {
    $overview.Cells.Item($r, 4).Value = "2016-20-12 10:20:47"
}

$detailRows = @(7, 10, 11, 12, 13, 14, 15, 16)
foreach ($r in $detailRows) {
    $zhcn.Cells.Item($r, 5).Value = "2016-03-12 10:20:42"
}
foreach ($r in $detailRows) {
    $dede.Cells.Item($r, 5).Value = "2016-03-12 10:20:47"
}
